# Manual_gating.xlsx — "Test analysis without single colors"
#
# Summary of the change:
#   1. Duplicate the original "Final" sheet (before any edits) into a new
#      sheet named "with_SC", inserted right after "Final" and before
#      "Sheet1".
#   2. On the "Final" sheet: swap a handful of manual-count values between
#      population rows (C4<->C5, C8<->C9, C11<->C12) and remove the Pop15
#      row (row 16) entirely, leaving the totals row at its original row
#      number (17).
#   3. Restore the expected active-cell selections on both sheets.

$wb = $excel.ActiveWorkbook

$final = $wb.Worksheets.Item("Final")

# --- Step 1: snapshot the untouched "Final" sheet as "with_SC" -------------
# Copy() places the duplicate immediately after the source sheet, i.e.
# between "Final" and "Sheet1" - exactly the order we need.
$final.Copy([System.Reflection.Missing]::Value, $final)
$withSC = $wb.Worksheets.Item("Final (2)")
$withSC.Name = "with_SC"

# Match the recorded (non-active) selection on the new sheet: a full
# A1:E17 range selection.
$withSC.Activate()
$withSC.Range("A1:E17").Select()

# --- Step 2: edit the live "Final" sheet ------------------------------------
$final.Activate()

# Swap manual counts for Pop2 / Pop3 (rows 4 and 5).
$final.Range("C4").Value2 = 58
$final.Range("C5").Value2 = 89

# Swap manual counts for Pop6 / Pop7 (rows 8 and 9).
$final.Range("C8").Value2 = 0
$final.Range("C9").Value2 = 4

# Swap manual counts for Pop9 / Pop11 (rows 11 and 12).
$final.Range("C11").Value2 = 2517
$final.Range("C12").Value2 = 295

# Drop the Pop15 row (row 16) completely - clear it instead of a
# shift-delete so the totals row below keeps its original row number (17).
$final.Range("A16:E16").ClearContents()

# --- Step 3: restore the active selection on "Final" ------------------------
$final.Range("E1").Select()
